$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A7").Value = "sprite animation script"
$ws.Range("A8").Value = "spawn scripts"
$ws.Range("A9").Value = "set board space to noone when dead"

$ws.Activate()
$ws.Range("A9").Select()
